$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new portfolio data row (row 20) for 2025-09-04.
# Force column A to be stored as text (matching the existing date-as-text
# column) rather than letting Excel auto-convert the string to a date serial.
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "2025-09-04"
$ws.Range("B20").Value = 57.41999816894531
$ws.Range("C20").Value = 687.5
$ws.Range("D20").Value = 326.25
